# SubSequentRenewal46_UT_SS.xlsx refresh
#
# The VIN upload test-data sheet is updated so the "symbol" table year moves
# from 2000 to 2017 and the previously-generic single-letter symbol codes
# (which happened to collide with other codes already in the workbook) are
# replaced with unique per-row BI/PD/UM/MP codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (VERSION) - symbol table year: SYMBOL_2000 -> SYMBOL_2017
$ws.Range("B2").Value = "SYMBOL_2017"
$ws.Range("B3").Value = "SYMBOL_2017"
$ws.Range("B4").Value = "SYMBOL_2017"
$ws.Range("B5").Value = "SYMBOL_2017"

# Columns AC:AF (BI_SYMBOL, PD_SYMBOL, UM_SYMBOL, MP_SYMBOL) - give each row
# its own distinct symbol code instead of reusing a single shared letter code.
$ws.Range("AC2").Value = "BI001"
$ws.Range("AD2").Value = "PD001"
$ws.Range("AE2").Value = "UM001"
$ws.Range("AF2").Value = "MP001"

$ws.Range("AC3").Value = "BI002"
$ws.Range("AD3").Value = "PD002"
$ws.Range("AE3").Value = "UM002"
$ws.Range("AF3").Value = "MP002"

$ws.Range("AC4").Value = "BI003"
$ws.Range("AD4").Value = "PD003"
$ws.Range("AE4").Value = "UM003"
$ws.Range("AF4").Value = "MP003"

$ws.Range("AC5").Value = "BI004"
$ws.Range("AD5").Value = "PD004"
$ws.Range("AE5").Value = "UM004"
$ws.Range("AF5").Value = "MP004"

# Leave the cursor on B9, matching the saved selection in the workbook.
$ws.Range("B9").Select()
